$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 65.06860185724791
$ws.Range("C2").Value = 94.66552152454085
$ws.Range("D2").Value = 98.82686518400283
$ws.Range("E2").Value = 98.28324009918686
$ws.Range("F2").Value = 98.42666309375164
$ws.Range("G2").Value = 97.57916325379115
$ws.Range("H2").Value = 96.13041750045974

$ws.Range("B3").Value = 67.2061853182154
$ws.Range("C3").Value = 94.16011017150215
$ws.Range("D3").Value = 99.82767734311415
$ws.Range("E3").Value = 98.91649001547911
$ws.Range("F3").Value = 98.63142936162016
$ws.Range("G3").Value = 97.65793611794369
$ws.Range("H3").Value = 96.1895233645517

$ws.Range("B4").Value = 86.40063628940017
$ws.Range("C4").Value = 93.08738621875607
$ws.Range("D4").Value = 98.78468112208979
$ws.Range("E4").Value = 98.78725109838855
$ws.Range("F4").Value = 98.4703883880736
$ws.Range("G4").Value = 97.53175931721469
$ws.Range("H4").Value = 96.11969775252989

$ws.Range("B5").Value = 76.1140395382857
$ws.Range("C5").Value = 92.38833927082841
$ws.Range("D5").Value = 98.760457854393
$ws.Range("E5").Value = 98.88946130907308
$ws.Range("F5").Value = 98.41562240624651
$ws.Range("G5").Value = 97.47559909487194
$ws.Range("H5").Value = 96.0760260454833

$ws.Range("B6").Value = 75.35715950445983
$ws.Range("C6").Value = 95.21512593686357
$ws.Range("D6").Value = 98.70946122742986
$ws.Range("E6").Value = 98.90334474522793
$ws.Range("F6").Value = 98.39272527353231
$ws.Range("G6").Value = 97.56543269061787
$ws.Range("H6").Value = 96.12280794758816
